# Apply cryptocurrency price/volume list refresh (GitHub Actions scheduled update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.926.19"
$ws.Range("E2").Value = "  -3.37%  "
$ws.Range("D3").Value = "3.359.26"
$ws.Range("E3").Value = "  -2.87%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "566.09"
$ws.Range("E5").Value = "  -2.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.59"
$ws.Range("E6").Value = "  -0.33%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  +1.10%  "
$ws.Range("E10").Value = "  -1.74%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.413"
$ws.Range("E11").Value = "  +1.14%  "
$ws.Range("D12").Value = "3.935.81"
$ws.Range("E12").Value = "  -2.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.02"
$ws.Range("E14").Value = "  -1.75%  "
$ws.Range("D15").Value = "3.354.11"
$ws.Range("E15").Value = "  -2.85%  "
$ws.Range("D17").Value = "60.995.37"
$ws.Range("E17").Value = "  -3.31%  "
$ws.Range("E18").Value = "  -2.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.21"
$ws.Range("E19").Value = "  -3.05%  "
$ws.Range("E20").Value = "  -4.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "374.52"
$ws.Range("E21").Value = "  -3.92%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "75.24"
$ws.Range("E22").Value = "  +0.69%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.559"
$ws.Range("E23").Value = "  -0.79%  "
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").Value = "3.499.99"
$ws.Range("E25").Value = "  -2.64%  "
$ws.Range("E26").Value = "  -6.81%  "
$ws.Range("E27").Value = "  -4.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.998"
$ws.Range("E28").Value = "  -0.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.39"
$ws.Range("E29").Value = "  -4.41%  "
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("E31").Value = "  -1.89%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.70"
$ws.Range("E32").Value = "  -4.79%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "22.84"
$ws.Range("E33").Value = "  -2.33%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.30"
$ws.Range("E34").Value = "  -4.82%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.36"
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "169.88"
$ws.Range("E36").Value = "  -0.26%  "
$ws.Range("E37").Value = "  -5.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.78"
$ws.Range("E38").Value = "  -4.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "29.08"
$ws.Range("E39").Value = "  -8.81%  "
$ws.Range("D40").Value = "3.394.32"
$ws.Range("E40").Value = "  -2.89%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0744"
$ws.Range("E41").Value = "  -4.37%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.29"
$ws.Range("E42").Value = "  -1.45%  "
$ws.Range("E43").Value = "  -4.29%  "
$ws.Range("E45").Value = "  -3.47%  "
$ws.Range("E46").Value = "  -6.28%  "
$ws.Range("D47").Value = "2.485.75"
$ws.Range("E48").Value = "  -3.34%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.50"
$ws.Range("E49").Value = "  -1.08%  "
$ws.Range("B50").Value = "FirstDigitalUSD"
$ws.Range("C50").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  -0.03%  "
$ws.Range("E51").Value = "  -2.61%  "

Write-Host "Applied crypto list update"
